$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scores = @{
    2  = 250
    4  = 550
    5  = 800
    6  = 800
    7  = 500
    8  = 600
    9  = 800
    10 = 800
    11 = 457.14
    13 = 600
    14 = 800
    15 = 800
    18 = 600
    19 = 800
    20 = 800
    21 = 600
    22 = 450
    23 = 685.71
}

foreach ($row in $scores.Keys) {
    $ws.Range("D$row").Value = $scores[$row]
}
